$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '30.708.83'
$ws.Cells.Item(2, 5).Value = '  +2.48%  '

$ws.Cells.Item(3, 4).Value = '1.889.21'
$ws.Cells.Item(3, 5).Value = '  +0.68%  '

$ws.Cells.Item(4, 5).Value = '  +0.15%  '

$ws.Cells.Item(5, 4).Formula = "'247.37"
$ws.Cells.Item(5, 5).Value = '  +2.16%  '

$ws.Cells.Item(6, 4).Formula = "'1.000"
$ws.Cells.Item(6, 5).Value = '  +0.06%  '

$ws.Cells.Item(7, 4).Formula = "'0.4943"
$ws.Cells.Item(7, 5).Value = '  +0.32%  '

$ws.Cells.Item(8, 4).Formula = "'0.2959"
$ws.Cells.Item(8, 5).Value = '  +1.53%  '

$ws.Cells.Item(9, 4).Formula = "'0.06823"
$ws.Cells.Item(9, 5).Value = '  +2.90%  '

$ws.Cells.Item(10, 4).Value = '1.887.58'
$ws.Cells.Item(10, 5).Value = '  +0.59%  '

$ws.Cells.Item(11, 4).Formula = "'17.21"
$ws.Cells.Item(11, 5).Value = '  +2.99%  '

$ws.Cells.Item(12, 4).Formula = "'0.07240"

$ws.Cells.Item(13, 4).Formula = "'91.71"
$ws.Cells.Item(13, 5).Value = '  +6.40%  '

$ws.Cells.Item(14, 4).Formula = "'5.079"
$ws.Cells.Item(14, 5).Value = '  +3.87%  '

$ws.Cells.Item(15, 4).Formula = "'0.6785"
$ws.Cells.Item(15, 5).Value = '  +2.00%  '

$ws.Cells.Item(16, 4).Value = '30.674.79'
$ws.Cells.Item(16, 5).Value = '  +2.48%  '

$ws.Cells.Item(17, 4).Formula = "'0.000007988"
$ws.Cells.Item(17, 5).Value = '  +1.80%  '

$ws.Cells.Item(18, 5).Value = '  +0.18%  '

$ws.Cells.Item(19, 4).Formula = "'13.23"

$ws.Cells.Item(20, 4).Value = '2.133.70'
$ws.Cells.Item(20, 5).Value = '  +0.57%  '

$ws.Cells.Item(21, 5).Value = '  +0.29%  '

$ws.Cells.Item(22, 4).Formula = "'4.831"
$ws.Cells.Item(22, 5).Value = '  +1.30%  '

$ws.Cells.Item(23, 4).Formula = "'188.82"
$ws.Cells.Item(23, 5).Value = '  +33.93%  '

$ws.Cells.Item(24, 4).Formula = "'6.060"
$ws.Cells.Item(24, 5).Value = '  +5.38%  '

$ws.Cells.Item(25, 4).Formula = "'9.355"
$ws.Cells.Item(25, 5).Value = '  +3.42%  '

$ws.Cells.Item(26, 4).Formula = "'156.34"
$ws.Cells.Item(26, 5).Value = '  +4.15%  '

$ws.Cells.Item(27, 4).Formula = "'19.17"
$ws.Cells.Item(27, 5).Value = '  +12.86%  '

$ws.Cells.Item(28, 5).Value = '  -0.18%  '

$ws.Cells.Item(29, 4).Formula = "'1.398"
$ws.Cells.Item(29, 5).Value = '  +0.33%  '

$ws.Cells.Item(30, 4).Formula = "'4.304"
$ws.Cells.Item(30, 5).Value = '  +2.86%  '

$ws.Cells.Item(31, 4).Formula = "'0.09012"
$ws.Cells.Item(31, 5).Value = '  +3.25%  '

$ws.Cells.Item(32, 4).Formula = "'4.015"
$ws.Cells.Item(32, 5).Value = '  +1.54%  '

$ws.Cells.Item(33, 4).Formula = "'0.05189"

$ws.Cells.Item(34, 4).Formula = "'0.7441"

$ws.Cells.Item(35, 5).Value = '  +0.25%  '

$ws.Cells.Item(36, 4).Formula = "'2.732"
$ws.Cells.Item(36, 5).Value = '  +2.31%  '

$ws.Cells.Item(37, 4).Formula = "'0.01840"
$ws.Cells.Item(37, 5).Value = '  +3.01%  '

$ws.Cells.Item(38, 4).Formula = "'2.669"
$ws.Cells.Item(38, 5).Value = '  -0.50%  '

$ws.Cells.Item(39, 4).Formula = "'2.157"
$ws.Cells.Item(39, 5).Value = '  -0.57%  '

$ws.Cells.Item(40, 4).Formula = "'0.9384"
$ws.Cells.Item(40, 5).Value = '  +0.98%  '

$ws.Cells.Item(41, 5).Value = '  +4.45%  '

$ws.Cells.Item(42, 4).Formula = "'105.33"
$ws.Cells.Item(42, 5).Value = '  +2.58%  '

$ws.Cells.Item(43, 2).Value = 'FraxShare'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(43, 4).Formula = "'5.776"
$ws.Cells.Item(43, 5).Value = '  +0.21%  '

$ws.Cells.Item(44, 2).Value = 'PaxDollar'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(44, 4).Formula = "'1.001"
$ws.Cells.Item(44, 5).Value = '  +0.17%  '

$ws.Cells.Item(45, 4).Formula = "'7.636"
$ws.Cells.Item(45, 5).Value = '  +2.90%  '

$ws.Cells.Item(46, 4).Formula = "'0.1342"
$ws.Cells.Item(46, 5).Value = '  +5.95%  '

$ws.Cells.Item(47, 5).Value = '  +3.19%  '

$ws.Cells.Item(48, 2).Value = 'EnergySwap'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(48, 4).Formula = "'8.701"
$ws.Cells.Item(48, 5).Value = '  +4.86%  '

$ws.Cells.Item(49, 2).Value = 'NEARProtocol'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(49, 4).Formula = "'1.429"
$ws.Cells.Item(49, 5).Value = '  +7.26%  '

$ws.Cells.Item(50, 4).Formula = "'0.3941"
$ws.Cells.Item(50, 5).Value = '  +4.48%  '

$ws.Cells.Item(51, 4).Formula = "'33.54"
$ws.Cells.Item(51, 5).Value = '  +3.28%  '
